# "All test to run changes in testng"
# Updates the Jira-tracking sheet with the current set of Key/Error rows
# and refreshes the random per-run id used on the Admin sheet.

$wb = $excel.ActiveWorkbook

# --- Admin sheet: the generated run id (Username column) changes ----------
$wsAdmin = $wb.Worksheets.Item("Admin")
$wsAdmin.Range("D2").Value = "8448782A"

# --- Jira sheet: rows for Key / Error description --------------------------
$wsJira = $wb.Worksheets.Item("Jira")

# Row 3: PersonalDetails -> Leave_Accept
$wsJira.Range("A3").Value = "Leave_Accept"

# Row 4: Recruitment_HiredList -> Recruitment_Rejected1, and the
# "Error description" cell now holds an explicit (blank) text value
# instead of being completely empty.
$wsJira.Range("A4").Value = "Recruitment_Rejected1"
$wsJira.Range("B4").Value = "'"
$wsJira.Range("B4").Style = "Normal"
